# Logic of adding extra amount and extra hours from excel sheet to database.
# Update existing check-in/out rows and append new attendance rows (3-9)
# mirroring the shape of row 2 (extra hour entries).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 adjustments
$ws.Range("A1").Value = 1
$ws.Range("C1").Value = 44542.666666666664

# Row 2 adjustment
$ws.Range("B2").Value = 44542.333333333336

# Append rows 3-9, cloning row 2's values/format (Copy preserves the cell
# style so the date/time number format carries over correctly).
for ($r = 3; $r -le 9; $r++) {
    $ws.Cells.Item(2, 1).Copy($ws.Cells.Item($r, 1))
    $ws.Cells.Item(2, 2).Copy($ws.Cells.Item($r, 2))
    $ws.Cells.Item(2, 3).Copy($ws.Cells.Item($r, 3))
    $ws.Cells.Item(2, 4).Copy($ws.Cells.Item($r, 4))
    $ws.Cells.Item(2, 5).Copy($ws.Cells.Item($r, 5))
}

# Match the saved selection/active cell from the workbook.
$ws.Range("B1").Select()
